$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to be treated as text so that
# numeric-looking strings (e.g. "0.991", "211.56") are not silently
# converted to floating point numbers with rounding artifacts.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "27.516.76"
$ws.Range("E2").Value = "  +1.69%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.573.64"
$ws.Range("E3").Value = "  +0.42%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "0.991"
$ws.Range("E4").Value = "  -1.41%  "

# Row 5 - BNB
$ws.Range("D5").Value = "211.56"
$ws.Range("E5").Value = "  +1.47%  "

# Row 6 - XRP
$ws.Range("D6").Value = "0.493"
$ws.Range("E6").Value = "  +0.28%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -1.47%  "

# Row 8 - Solana
$ws.Range("D8").Value = "22.88"
$ws.Range("E8").Value = "  +3.53%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +0.72%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +0.31%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +1.33%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.796.03"
$ws.Range("E12").Value = "  +0.30%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.571.83"
$ws.Range("E13").Value = "  +0.00%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  -0.44%  "

# Row 15 - Polygon
$ws.Range("D15").Value = "0.520"
$ws.Range("E15").Value = "  -0.04%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "27.489.39"
$ws.Range("E16").Value = "  +1.53%  "

# Row 17 - Litecoin
$ws.Range("D17").Value = "62.50"
$ws.Range("E17").Value = "  +1.01%  "

# Row 18 - BitcoinCash
$ws.Range("D18").Value = "226.91"
$ws.Range("E18").Value = "  +5.15%  "

# Row 19 - Chainlink
$ws.Range("E19").Value = "  +1.31%  "

# Row 20 - ShibaInu
$ws.Range("E20").Value = "  +0.27%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  -1.43%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  -0.58%  "

# Row 23 - Avalanche
$ws.Range("E23").Value = "  +2.43%  "

# Row 24 - Toncoin
$ws.Range("D24").Value = "1.95"
$ws.Range("E24").Value = "  +0.18%  "

# Row 25 - Monero
$ws.Range("D25").Value = "150.71"
$ws.Range("E25").Value = "  -2.14%  "

# Row 26 - EthereumClassic
$ws.Range("D26").Value = "15.19"
$ws.Range("E26").Value = "  +0.99%  "

# Row 27 - Cosmos
$ws.Range("E27").Value = "  -0.30%  "

# Row 28 - Stellar
$ws.Range("E28").Value = "  +1.72%  "

# Row 29 - BinanceUSD
$ws.Range("D29").Value = "0.992"
$ws.Range("E29").Value = "  -1.31%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  +1.17%  "

# Row 31 - Hedera
$ws.Range("D31").Value = "0.0473"
$ws.Range("E31").Value = "  -0.38%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  +0.53%  "

# Row 33 - Maker
$ws.Range("D33").Value = "1.457.34"
$ws.Range("E33").Value = "  +2.43%  "

# Row 34 - InternetComputer(DFINITY)
$ws.Range("D34").Value = "3.15"
$ws.Range("E34").Value = "  -1.81%  "

# Row 35 - TrustWalletToken
$ws.Range("E35").Value = "  +3.48%  "

# Row 36 - LidoDAOToken
$ws.Range("E36").Value = "  -0.18%  "

# Row 37 - HuobiToken
$ws.Range("E37").Value = "  -0.85%  "

# Row 38 - VeChain
$ws.Range("E38").Value = "  +0.43%  "

# Row 39 - ImmutableX
$ws.Range("E39").Value = "  +1.37%  "

# Row 40 - ARBITRUM
$ws.Range("E40").Value = "  +0.26%  "

# Row 41 - MXToken
$ws.Range("E41").Value = "  -0.74%  "

# Row 42 - was FraxShare, now PaxDollar (rows 42/43 swapped with updated values)
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").Value = "0.992"
$ws.Range("E42").Value = "  -1.43%  "

# Row 43 - was PaxDollar, now FraxShare
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "5.65"
$ws.Range("E43").Value = "  -3.08%  "

# Row 44 - RenderToken
$ws.Range("E44").Value = "  +6.94%  "

# Row 45 - WEMIXToken
$ws.Range("D45").Value = "0.979"
$ws.Range("E45").Value = "  -2.40%  "

# Row 46 - Aave
$ws.Range("D46").Value = "64.33"
$ws.Range("E46").Value = "  -0.59%  "

# Row 47 - RocketPoolETH
$ws.Range("D47").Value = "1.708.20"
$ws.Range("E47").Value = "  +0.23%  "

# Row 48 - Quant
$ws.Range("D48").Value = "87.01"
$ws.Range("E48").Value = "  +0.13%  "

# Row 49 - BabyDogeCoin
$ws.Range("E49").Value = "  +0.78%  "

# Row 50 - Cronos
$ws.Range("E50").Value = "  +1.44%  "

# Row 51 - Algorand
$ws.Range("D51").Value = "0.0947"
$ws.Range("E51").Value = "  -1.74%  "

# Restore the original (default) cell style now that the text values
# have been written, so no stray number-format styling is left behind.
$ws.Range("D2:E51").Style = "Normal"
